$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1625
$ws.Range("D4").Value = 2036448
$ws.Range("D5").Value = 1905766
$ws.Range("D6").Value = 2115927
$ws.Range("D7").Value = 6058130

$ws.Range("E3:F7").ClearContents()
$ws.Range("D8:F9").ClearContents()

$ws.Range("F7").Select()
